# Weekly update: insert a new price record row for Betarraga
# (Feria Lagunitas de Puerto Montt) just above the existing row 469,
# pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 469 (shifts rows 469:505 -> 470:506)
$ws.Rows("469:469").Insert()

# Populate the new row with the latest weekly record
$ws.Cells.Item(469, 1).Value = 4
$ws.Cells.Item(469, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(469, 3).Value = "Los Lagos"
$ws.Cells.Item(469, 4).Value2 = 45106
$ws.Cells.Item(469, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(469, 5).Value = 10
$ws.Cells.Item(469, 6).Value = 100114014
$ws.Cells.Item(469, 7).Value = "Betarraga"
$ws.Cells.Item(469, 8).Value = "Sin especificar"
$ws.Cells.Item(469, 9).Value = "Primera"
$ws.Cells.Item(469, 10).Value = 500
$ws.Cells.Item(469, 11).Value = 1200
$ws.Cells.Item(469, 12).Value = 1200
$ws.Cells.Item(469, 13).Value = 1200
$ws.Cells.Item(469, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(469, 15).Value = "Región Metropolitana"
$ws.Cells.Item(469, 16).Value = 240
$ws.Cells.Item(469, 17).Value = 5
$ws.Cells.Item(469, 18).Value = "Hortaliza"
